# trading_journal.xlsx edit: "Add files via upload"
#
# Summary of the change (per the supplied OOXML diff):
#   - Sheet "Sheet2" (sheetId 5 / rId6) is renamed to "index tracking" and
#     becomes the active sheet (its tab is selected and its column C is
#     selected in its view); two new dated rows are appended to it, and
#     its column B is widened to fit the new text.
#   - A shared-string on Sheet1 (the trading journal "Reason for
#     loss/profit" note in E75) is corrected/extended, and a brand-new
#     trailing note is appended as row 76 in column E.
#   - The previously-active Sheet1 view loses its "active tab" status and
#     its cursor/scroll moves down near the newly added rows.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: fix the existing note and append a brand-new one ---------

$sheet1.Range("E75").Value = "I didn't notice the change in lot size so booked loss again always try to trade second trade after 12 PM keep in mind so that you will know market trend if early trade goes against your process"

$sheet1.Range("E76").Value = "In Trading risk management is the biggest and biggest capital and learning ever you will know till the end career so never ever forget"

# Leave the cursor near the new row, matching the author's final view.
$sheet1.Range("E76").Select()

# --- Sheet2: rename, add two new journal rows, widen column B ---------

$sheet2.Name = "index tracking"

# New dated rows - copy the date format from the row above (A5) so the
# new date cells reuse the existing short-date style instead of minting a
# new number format.
$sheet2.Range("A5").Copy()
$sheet2.Range("A6:A7").PasteSpecial(-4122)

$sheet2.Range("A6").Value = 45652
$sheet2.Range("B6").Value = "market was in a compressed zone so day after it there was a breakout"

$sheet2.Range("A7").Value = 45653
$sheet2.Range("B7").Value = "market was In a bit breakout in the morning session but again was consolidating in day time no move/wait and watch for next day "

$sheet2.Columns.Item(2).ColumnWidth = 108

# index tracking becomes the active sheet/tab, with column C selected.
$sheet2.Activate()
$sheet2.Range("C:C").Select()
